# Lecture time table edits:
# - Move "CD222-Dr. Wendy Stokes-Hall 5" from E4 (Tuesday, 2:00-3:30PM) to E7 (Friday, 2:00-3:30PM)
# - Move "CD222-sec-Hall 5" from B8 (Saturday, 9:00-10:30AM) to C3 (Monday, 10:30AM-12:00PM)
# - Move "CD222-lab-Hall 5" from D8 (Saturday, 12:30-2:00PM) to F5 (Wednesday, 3:30-5:00PM)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture existing values before clearing so moves are order-independent
$valE4 = $ws.Range("E4").Value2
$valB8 = $ws.Range("B8").Value2
$valD8 = $ws.Range("D8").Value2

# Clear the old locations
$ws.Range("E4").ClearContents()
$ws.Range("B8").ClearContents()
$ws.Range("D8").ClearContents()

# Write the values into their new locations
$ws.Range("E7").Value = $valE4
$ws.Range("C3").Value = $valB8
$ws.Range("F5").Value = $valD8
